# Update CDA Logical model for ST.r2b
# Edits the "Metadata" worksheet (Property/Value table):
#   - Version bumped
#   - Date bumped
#   - new "Jurisdiction" property row inserted after "Contact" and before "Description"

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# Update Version value (row 3, column B)
$ws.Cells.Item(3, 2).Value = "2.0.1-sd-202510-matchbox-patch"

# Update Date value (row 8, column B)
$ws.Cells.Item(8, 2).Value = "2025-10-29T22:15:57+01:00"

# Insert a new row for "Jurisdiction" right after the "Contact" row (row 10),
# pushing "Description" and everything below it down by one row.
$ws.Rows.Item(11).Insert()

# Match the formatting of the surrounding data rows (copy format from the row
# that is now directly below, i.e. the former "Description" row).
$ws.Range("A12:B12").Copy()
$ws.Range("A11:B11").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Cells.Item(11, 1).Value = "Jurisdiction"
$ws.Cells.Item(11, 2).Value = ""
